$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 24; $r -le 98; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
